# Auto-generated edit script for LOB1221.docx content reshuffle.
# Strategy: two-phase replace using unique placeholder tokens to avoid the
# Find/Replace engine mis-merging runs when a ReplaceWith text coincides with
# pre-existing text elsewhere in the document.
$d = $word.ActiveDocument

function Find-ReplaceInParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 1)
    if (-not $ok) {
        throw "Find/Replace failed in paragraph $paraIndex for: $oldText"
    }
}

# ---- Phase 1: stamp unique placeholders over every slot that changes ----
Find-ReplaceInParagraph 9 '5840671 - Francisco José Moreira Chaves' '@@PH_1_DOCENTE_TO_OBJPT@@'
Find-ReplaceInParagraph 6 'O objetivo da presente disciplina é introduzir os alunos no sistema normativo ambiental, conhecendo os princípios fundamentais do Direito Ambiental, sendo também capazes de analisar alguns dos instrumentos da Política Nacional de Meio Ambiente e discutir aspectos da legislação protetora dos recursos ambientais.' '@@PH_2_OBJPT_TO_RESUMOPT@@'
Find-ReplaceInParagraph 11 'Direitos ambiental constitucional; política nacional do meio ambiente' '@@PH_3_RESUMOPT_TO_PROGPT@@'
Find-ReplaceInParagraph 14 'Conceitos básicos; Princípios fundamentais de direito ambiental; evolução histórica da legislação ambiental brasileira; política nacional do meio ambiente ; código florestal brasileiro; política nacional de recursos hídricos; lei dos crimes ambientais; sistema nacional de unidades de conservação; tutela administrativa, civil e processual do meio ambiente; estudos de caso com aplicação da legislação ambiental vigente e necessária para o licenciamento de empreendimento em diversos estados brasileiros.' '@@PH_4_PROGPT_TO_METODO@@'
Find-ReplaceInParagraph 17 'As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático.' '@@PH_5_METODO_TO_CRITERIO@@'
Find-ReplaceInParagraph 17 'Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois.' '@@PH_6_CRITERIO_TO_NORMA@@'
Find-ReplaceInParagraph 17 'A nota final será composta pela média obtida da nota do período somada à nota de recuperação' '@@PH_7_NORMA_TO_BIBLIO@@'
Find-ReplaceInParagraph 19 'ALENZA G. J. F. Manual de Derecho Ambiental. Universidad Pública de Navarra, 2001. ^lANTUNES, P. B. Dano Ambiental: uma abordagem conceitual. Rio de Janeiro, Editora Lumen Juris, 2000.^lBUSTAMANTE A. J. Derecho Ambiental, Editorial Abeledo-perrot, Buenos Aires.^lCOSTA JR., P. J. Direito Penal Ecológico. Rio de Janeiro, Forense Universitária, 1996.^lCOSTA NETO, N. D. C.; BELLO FILHO, N. B.; e CASTRO E COSTA, F. D. Crimes e Infrações Administrativas Ambientais. Brasília: Brasília Jurídica, 2000. ^lDICIONÁRIO DE DIREITO AMBIENTAL - Terminologia das Leis do Meio Ambiente. Maria da Graça Krieger, Anna Maria Becker Maciel, João Carlos de Carvalho Rocha, Maria José Bocorny Finatto e Cleci Regina Bevilacqua. Editora Universidade/UFRGS.^lFIORILLO, C. A. P.; e RODRIGES, M. A. Manual de Direito Ambiental e legislação aplicável. São Paulo, Max Limonad, 1997.^lFIORILLO, C. A. P.; e RODRIGES, M. A. Direito Ambiental e Patrimônio Genético. Belo Horizonte: Del Rey, 1996.^lFREITAS, V. P. Águas - Aspectos Jurídicos e Ambientais. Curitiba, Juruá, 2000.^lFREITAS, V. P. Direito Administrativo e Meio Ambiente. Curitiba, Juruá, 1993.^lGOMES, C. L. S. P. Crimes Contra o Meio Ambiente: responsabilidade e sanção penal. 2ª edição, São Paulo, Editora Juarez de Oliveira, 1999.^lMACHADO, P. A. L. Direito Ambiental Brasileiro. 8ª Edição, Revista, atualizada e ampliada, São Paulo, Malheiros Editores, 2001.^lMORAES, A. Direito Constitucional.. 7ª ed. revista, ampliada e atualizada, com a EC n.º 24/99 - São Paulo, Atlas, 2000.' '@@PH_8_BIBLIO_TO_DOCENTE@@'
Find-ReplaceInParagraph 7 'The objective of this course is to introduce students to the environmental normative system, knowing the fundamental principles of Environmental Law, and also to analyze some of the instruments of the National Environmental Policy and discuss aspects of the protective legislation of environmental resources.' '@@PH_9_OBJEN_TO_RESUMOEN@@'
Find-ReplaceInParagraph 12 'Constitutional environmental law; National environment policy' '@@PH_10_RESUMOEN_TO_OBJEN@@'

# ---- Phase 2: replace each placeholder with its real final text ----
Find-ReplaceInParagraph 9 '@@PH_1_DOCENTE_TO_OBJPT@@' 'O objetivo da presente disciplina é introduzir os alunos no sistema normativo ambiental, conhecendo os princípios fundamentais do Direito Ambiental, sendo também capazes de analisar alguns dos instrumentos da Política Nacional de Meio Ambiente e discutir aspectos da legislação protetora dos recursos ambientais.'
Find-ReplaceInParagraph 6 '@@PH_2_OBJPT_TO_RESUMOPT@@' 'Direitos ambiental constitucional; política nacional do meio ambiente'
Find-ReplaceInParagraph 11 '@@PH_3_RESUMOPT_TO_PROGPT@@' 'Conceitos básicos; Princípios fundamentais de direito ambiental; evolução histórica da legislação ambiental brasileira; política nacional do meio ambiente ; código florestal brasileiro; política nacional de recursos hídricos; lei dos crimes ambientais; sistema nacional de unidades de conservação; tutela administrativa, civil e processual do meio ambiente; estudos de caso com aplicação da legislação ambiental vigente e necessária para o licenciamento de empreendimento em diversos estados brasileiros.'
Find-ReplaceInParagraph 14 '@@PH_4_PROGPT_TO_METODO@@' 'As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático.'
Find-ReplaceInParagraph 17 '@@PH_5_METODO_TO_CRITERIO@@' 'Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois.'
Find-ReplaceInParagraph 17 '@@PH_6_CRITERIO_TO_NORMA@@' 'A nota final será composta pela média obtida da nota do período somada à nota de recuperação'
Find-ReplaceInParagraph 17 '@@PH_7_NORMA_TO_BIBLIO@@' 'ALENZA G. J. F. Manual de Derecho Ambiental. Universidad Pública de Navarra, 2001. ^lANTUNES, P. B. Dano Ambiental: uma abordagem conceitual. Rio de Janeiro, Editora Lumen Juris, 2000.^lBUSTAMANTE A. J. Derecho Ambiental, Editorial Abeledo-perrot, Buenos Aires.^lCOSTA JR., P. J. Direito Penal Ecológico. Rio de Janeiro, Forense Universitária, 1996.^lCOSTA NETO, N. D. C.; BELLO FILHO, N. B.; e CASTRO E COSTA, F. D. Crimes e Infrações Administrativas Ambientais. Brasília: Brasília Jurídica, 2000. ^lDICIONÁRIO DE DIREITO AMBIENTAL - Terminologia das Leis do Meio Ambiente. Maria da Graça Krieger, Anna Maria Becker Maciel, João Carlos de Carvalho Rocha, Maria José Bocorny Finatto e Cleci Regina Bevilacqua. Editora Universidade/UFRGS.^lFIORILLO, C. A. P.; e RODRIGES, M. A. Manual de Direito Ambiental e legislação aplicável. São Paulo, Max Limonad, 1997.^lFIORILLO, C. A. P.; e RODRIGES, M. A. Direito Ambiental e Patrimônio Genético. Belo Horizonte: Del Rey, 1996.^lFREITAS, V. P. Águas - Aspectos Jurídicos e Ambientais. Curitiba, Juruá, 2000.^lFREITAS, V. P. Direito Administrativo e Meio Ambiente. Curitiba, Juruá, 1993.^lGOMES, C. L. S. P. Crimes Contra o Meio Ambiente: responsabilidade e sanção penal. 2ª edição, São Paulo, Editora Juarez de Oliveira, 1999.^lMACHADO, P. A. L. Direito Ambiental Brasileiro. 8ª Edição, Revista, atualizada e ampliada, São Paulo, Malheiros Editores, 2001.^lMORAES, A. Direito Constitucional.. 7ª ed. revista, ampliada e atualizada, com a EC n.º 24/99 - São Paulo, Atlas, 2000.'
Find-ReplaceInParagraph 19 '@@PH_8_BIBLIO_TO_DOCENTE@@' '5840671 - Francisco José Moreira Chaves'
Find-ReplaceInParagraph 7 '@@PH_9_OBJEN_TO_RESUMOEN@@' 'Constitutional environmental law; National environment policy'
Find-ReplaceInParagraph 12 '@@PH_10_RESUMOEN_TO_OBJEN@@' 'The objective of this course is to introduce students to the environmental normative system, knowing the fundamental principles of Environmental Law, and also to analyze some of the instruments of the National Environmental Policy and discuss aspects of the protective legislation of environmental resources.'

Write-Output "Done."
